# Auto-generated Excel COM-interop script
# Applies numeric corrections to currentAveragePrice / Leve price / profit columns
# across multiple crafting-class sheets (Sheets/Tiamat_Profits.xlsx export).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1667.909
$ws.Range("I4").Value = 905.44446
$ws.Range("J4").Value = 5099
$ws.Range("K4").Value = 905.44446
$ws.Range("L4").Value = 5099
$ws.Range("M4").Value = -791.44446
$ws.Range("N4").Value = -5327

$ws.Range("H69").Value = 2454210.2
$ws.Range("I69").Value = 14705882
$ws.Range("J69").Value = 3876
$ws.Range("K69").Value = 44117646
$ws.Range("L69").Value = 11628
$ws.Range("M69").Value = -44116772
$ws.Range("N69").Value = -13376

$ws.Range("H72").Value = 2454210.2
$ws.Range("I72").Value = 14705882
$ws.Range("J72").Value = 3876
$ws.Range("K72").Value = 132352938
$ws.Range("L72").Value = 34884
$ws.Range("M72").Value = -132348570
$ws.Range("N72").Value = -43620

$ws.Range("H76").Value = 17860320
$ws.Range("I76").Value = 23812444
$ws.Range("J76").Value = 3947
$ws.Range("K76").Value = 23812444
$ws.Range("L76").Value = 3947
$ws.Range("M76").Value = -23812129
$ws.Range("N76").Value = -4577

$ws.Range("H79").Value = 17860320
$ws.Range("I79").Value = 23812444
$ws.Range("J79").Value = 3947
$ws.Range("K79").Value = 23812444
$ws.Range("L79").Value = 3947
$ws.Range("M79").Value = -23811352
$ws.Range("N79").Value = -6131

$ws.Range("H80").Value = 4548154
$ws.Range("I80").Value = 2022.6666
$ws.Range("K80").Value = 6067.9998
$ws.Range("M80").Value = -5069.9998

$ws.Range("H82").Value = 2260.625
$ws.Range("I82").Value = 558.8
$ws.Range("J82").Value = 3034.182
$ws.Range("K82").Value = 1676.4
$ws.Range("L82").Value = 9102.545999999998
$ws.Range("M82").Value = -1270.4
$ws.Range("N82").Value = -9914.545999999998

$ws.Range("H83").Value = 4548154
$ws.Range("I83").Value = 2022.6666
$ws.Range("K83").Value = 18203.9994
$ws.Range("M83").Value = -13211.9994

$ws.Range("H85").Value = 2260.625
$ws.Range("I85").Value = 558.8
$ws.Range("J85").Value = 3034.182
$ws.Range("K85").Value = 1676.4
$ws.Range("L85").Value = 9102.545999999998
$ws.Range("M85").Value = -272.3999999999999
$ws.Range("N85").Value = -11910.546

$ws.Range("H88").Value = 7000
$ws.Range("I88").Value = 6933.3335
$ws.Range("J88").Value = 7066.6665
$ws.Range("K88").Value = 6933.3335
$ws.Range("L88").Value = 7066.6665
$ws.Range("M88").Value = -6527.3335
$ws.Range("N88").Value = -7878.6665

$ws.Range("H91").Value = 7000
$ws.Range("I91").Value = 6933.3335
$ws.Range("J91").Value = 7066.6665
$ws.Range("K91").Value = 6933.3335
$ws.Range("L91").Value = 7066.6665
$ws.Range("M91").Value = -5529.3335
$ws.Range("N91").Value = -9874.666499999999

$ws.Range("H132").Value = 373549.6
$ws.Range("I132").Value = 3393.08
$ws.Range("J132").Value = 5000506
$ws.Range("K132").Value = 10179.24
$ws.Range("L132").Value = 15001518
$ws.Range("M132").Value = -7649.24
$ws.Range("N132").Value = -15006578

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

$ws.Range("H32").Value = 171702.28
$ws.Range("I32").Value = 167414.9
$ws.Range("K32").Value = 167414.9
$ws.Range("M32").Value = -167127.9

$ws.Range("H45").Value = 1487.091
$ws.Range("I45").Value = 1479.0526
$ws.Range("J45").Value = 1538
$ws.Range("K45").Value = 1479.0526
$ws.Range("L45").Value = 1538
$ws.Range("M45").Value = -1102.0526
$ws.Range("N45").Value = -2292

$ws.Range("H76").Value = 59463.5
$ws.Range("J76").Value = 59463.5
$ws.Range("L76").Value = 59463.5
$ws.Range("N76").Value = -60139.5

$ws.Range("H79").Value = 59463.5
$ws.Range("J79").Value = 59463.5
$ws.Range("L79").Value = 59463.5
$ws.Range("N79").Value = -61803.5

$ws.Range("H122").Value = 1478.8
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 1473.5
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 4420.5
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -9320.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws.Range("H105").Value = 996826.25
$ws.Range("I105").Value = 1328668.4
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1328668.4
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = -1326921.4
$ws.Range("N105").Value = -4794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1020
$ws.Range("I7").Value = 1542.8572
$ws.Range("J7").Value = 105
$ws.Range("K7").Value = 1542.8572
$ws.Range("L7").Value = 105
$ws.Range("M7").Value = -1429.8572
$ws.Range("N7").Value = -331

$ws.Range("H51").Value = 8351.552
$ws.Range("J51").Value = 8351.552
$ws.Range("L51").Value = 8351.552
$ws.Range("N51").Value = -9823.552

$ws.Range("H61").Value = 8351.552
$ws.Range("J61").Value = 8351.552
$ws.Range("L61").Value = 8351.552
$ws.Range("N61").Value = -9047.552

$ws.Range("H74").Value = 11992.353
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 12429.375
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 12429.375
$ws.Range("M74").Value = -4126
$ws.Range("N74").Value = -14177.375

$ws.Range("H77").Value = 11992.353
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 12429.375
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 37288.125
$ws.Range("M77").Value = -10632
$ws.Range("N77").Value = -46024.125

$ws.Range("H88").Value = 41794
$ws.Range("J88").Value = 41794
$ws.Range("L88").Value = 41794
$ws.Range("N88").Value = -42606

$ws.Range("H91").Value = 41794
$ws.Range("J91").Value = 41794
$ws.Range("L91").Value = 41794
$ws.Range("N91").Value = -44602

$ws.Range("H134").Value = 16130971
$ws.Range("I134").Value = 1694.9524
$ws.Range("J134").Value = 50002452
$ws.Range("K134").Value = 5084.857199999999
$ws.Range("L134").Value = 150007356
$ws.Range("M134").Value = -2549.857199999999
$ws.Range("N134").Value = -150012426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4229617
$ws.Range("I4").Value = 37002200
$ws.Range("J4").Value = 896.7742
$ws.Range("K4").Value = 111006600
$ws.Range("L4").Value = 2690.3226
$ws.Range("M4").Value = -111006488
$ws.Range("N4").Value = -2914.3226

$ws.Range("H61").Value = 1100
$ws.Range("I61").Value = 100
$ws.Range("J61").Value = 1300
$ws.Range("K61").Value = 300
$ws.Range("L61").Value = 3900
$ws.Range("M61").Value = -85
$ws.Range("N61").Value = -4330

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 20000
$ws.Range("J5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("N5").Value = -20224

$ws.Range("H102").Value = 81105.2
$ws.Range("I102").Value = 26378
$ws.Range("J102").Value = 300014
$ws.Range("K102").Value = 26378
$ws.Range("L102").Value = 300014
$ws.Range("M102").Value = -24756
$ws.Range("N102").Value = -303258

$ws.Range("H122").Value = 4425.2085
$ws.Range("I122").Value = 4066.6667
$ws.Range("J122").Value = 5022.778
$ws.Range("K122").Value = 12200.0001
$ws.Range("L122").Value = 15068.334
$ws.Range("M122").Value = -9750.000100000001
$ws.Range("N122").Value = -19968.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 70002
$ws.Range("J2").Value = 70002
$ws.Range("L2").Value = 70002
$ws.Range("N2").Value = -70226

$ws.Range("H64").Value = 19998.572
$ws.Range("J64").Value = 19998.572
$ws.Range("L64").Value = 19998.572
$ws.Range("N64").Value = -20448.572

$ws.Range("H67").Value = 19998.572
$ws.Range("J67").Value = 19998.572
$ws.Range("L67").Value = 19998.572
$ws.Range("N67").Value = -21558.572

$ws.Range("H122").Value = 2853.3333
$ws.Range("I122").Value = 2280
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 6840
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -4390
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 356031
$ws.Range("I132").Value = 89548.30499999999
$ws.Range("J132").Value = 913222.0600000001
$ws.Range("K132").Value = 268644.915
$ws.Range("L132").Value = 2739666.18
$ws.Range("M132").Value = -266114.915
$ws.Range("N132").Value = -2744726.18

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4000500
$ws.Range("J2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("N2").Value = -1224

$ws.Range("H81").Value = 2008.0834
$ws.Range("I81").Value = 1440.75
$ws.Range("J81").Value = 2575.4167
$ws.Range("K81").Value = 2881.5
$ws.Range("L81").Value = 5150.8334
$ws.Range("M81").Value = -1820.5
$ws.Range("N81").Value = -7272.8334

$ws.Range("H82").Value = 49963.668
$ws.Range("J82").Value = 49963.668
$ws.Range("L82").Value = 49963.668
$ws.Range("N82").Value = -50729.668

$ws.Range("H84").Value = 2008.0834
$ws.Range("I84").Value = 1440.75
$ws.Range("J84").Value = 2575.4167
$ws.Range("K84").Value = 14407.5
$ws.Range("L84").Value = 25754.167
$ws.Range("M84").Value = -9103.5
$ws.Range("N84").Value = -36362.167

$ws.Range("H85").Value = 49963.668
$ws.Range("J85").Value = 49963.668
$ws.Range("L85").Value = 49963.668
$ws.Range("N85").Value = -52615.668

$ws.Range("H122").Value = 4363.769
$ws.Range("I122").Value = 3698.439
$ws.Range("J122").Value = 6843.636
$ws.Range("K122").Value = 11095.317
$ws.Range("L122").Value = 20530.908
$ws.Range("M122").Value = -8645.316999999999
$ws.Range("N122").Value = -25430.908

$ws.Range("H132").Value = 5014.931
$ws.Range("I132").Value = 1330.1666
$ws.Range("J132").Value = 22701.8
$ws.Range("K132").Value = 3990.4998
$ws.Range("L132").Value = 68105.39999999999
$ws.Range("M132").Value = -1460.4998
$ws.Range("N132").Value = -73165.39999999999

